$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-looking numbers as TEXT in the source
# sheet (e.g. "299.52"), same cell type as the thousands-grouped values like
# "42.377.41" that Excel could never parse as a number anyway. Force those
# particular cells to the Text format first so assigning the new reading
# doesn't get auto-coerced into a numeric value/type by Excel.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('D2').Value = '42.377.41'
$ws.Range('E2').Value = '  -2.90%  '
$ws.Range('D3').Value = '2.270.31'
$ws.Range('E3').Value = '  -4.77%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '299.52'
$ws.Range('E5').Value = '  -3.35%  '
$ws.Range('D6').Value = '96.32'
$ws.Range('E6').Value = '  -7.85%  '
$ws.Range('D7').Value = '0.502'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  -4.65%  '
$ws.Range('D10').Value = '33.71'
$ws.Range('E10').Value = '  -6.55%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '50.63'
$ws.Range('E11').Value = '  -5.24%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '0.0784'
$ws.Range('E12').Value = '  -3.73%  '
$ws.Range('D13').Value = '0.112'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '6.64'
$ws.Range('E14').Value = '  -5.12%  '
$ws.Range('D15').Value = '2.622.01'
$ws.Range('E15').Value = '  -4.71%  '
$ws.Range('D16').Value = '15.21'
$ws.Range('E16').Value = '  -2.88%  '
$ws.Range('D17').Value = '2.267.58'
$ws.Range('E17').Value = '  -4.62%  '
$ws.Range('D18').Value = '0.782'
$ws.Range('E18').Value = '  -3.65%  '
$ws.Range('D19').Value = '42.308.01'
$ws.Range('E19').Value = '  -2.98%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = '11.40'
$ws.Range('E21').Value = '  -4.23%  '
$ws.Range('D22').Value = '5.97'
$ws.Range('E22').Value = '  -5.61%  '
$ws.Range('D23').Value = '66.47'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '233.67'
$ws.Range('E24').Value = '  -3.18%  '
$ws.Range('D25').Value = '1.93'
$ws.Range('E25').Value = '  -5.98%  '
$ws.Range('D26').Value = '2.48'
$ws.Range('E26').Value = '  -5.38%  '
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').Value = '24.34'
$ws.Range('E28').Value = '  -5.87%  '
$ws.Range('E29').Value = '  +3.04%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '164.09'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '33.71'
$ws.Range('E31').Value = '  -7.85%  '
$ws.Range('D32').Value = '9.05'
$ws.Range('E32').Value = '  -5.18%  '
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('D34').Value = '4.93'
$ws.Range('E34').Value = '  -6.33%  '
$ws.Range('E35').Value = '  -4.72%  '
$ws.Range('D36').Value = '0.0694'
$ws.Range('E36').Value = '  -6.26%  '
$ws.Range('D37').Value = '4.36'
$ws.Range('E37').Value = '  -6.48%  '
$ws.Range('D38').Value = '2.81'
$ws.Range('E38').Value = '  -10.14%  '
$ws.Range('D39').Value = '16.04'
$ws.Range('E39').Value = '  -12.50%  '
$ws.Range('D40').Value = '0.0999'
$ws.Range('E40').Value = '  -5.50%  '
$ws.Range('D41').Value = '1.76'
$ws.Range('E41').Value = '  -9.45%  '
$ws.Range('D42').Value = '0.109'
$ws.Range('E42').Value = '  -3.92%  '
$ws.Range('D43').Value = '2.41'
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('D44').Value = '1.960.50'
$ws.Range('E44').Value = '  -3.68%  '
$ws.Range('D45').Value = '0.0281'
$ws.Range('E45').Value = '  -3.32%  '
$ws.Range('D46').Value = '17.70'
$ws.Range('E46').Value = '  -9.82%  '
$ws.Range('D47').Value = '9.68'
$ws.Range('E47').Value = '  -8.55%  '
$ws.Range('D48').Value = '2.81'
$ws.Range('E48').Value = '  -10.38%  '
$ws.Range('D49').Value = '2.84'
$ws.Range('D50').Value = '4.67'
$ws.Range('E50').Value = '  -1.72%  '
$ws.Range('D51').Value = '2.502.53'
$ws.Range('E51').Value = '  -3.99%  '
